$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97
$ws.Cells.Item(97, 2).Value = 6809780
$ws.Cells.Item(97, 6).Value = 'WaaslandBeveren'
$ws.Cells.Item(97, 7).Value = 'FCV Dender EH'
$ws.Cells.Item(97, 8).Value = 3
$ws.Cells.Item(97, 10).Value = 'H'
$ws.Cells.Item(97, 11).Value = 1.95
$ws.Cells.Item(97, 12).Value = 3.6
$ws.Cells.Item(97, 13).Value = 3.3
$ws.Cells.Item(97, 14).Value = 2.15
$ws.Cells.Item(97, 16).Value = 3
$ws.Cells.Item(97, 17).Value = -0.25
$ws.Cells.Item(97, 18).Value = 1.9
$ws.Cells.Item(97, 19).Value = 1.9
$ws.Cells.Item(97, 20).Value = 2.75
$ws.Cells.Item(97, 21).Value = 1.975
$ws.Cells.Item(97, 22).Value = 1.825
$ws.Cells.Item(97, 23).Value = 1.15
$ws.Cells.Item(97, 25).Value = -1
$ws.Cells.Item(97, 26).Value = 0.8999999999999999
$ws.Cells.Item(97, 27).Value = -1
$ws.Cells.Item(97, 28).Value = 0.9750000000000001
# Row 98
$ws.Cells.Item(98, 2).Value = 6809783
$ws.Cells.Item(98, 6).Value = 'Lommel'
$ws.Cells.Item(98, 7).Value = 'Patro Eisden Maasmechelen'
$ws.Cells.Item(98, 8).Value = 1
$ws.Cells.Item(98, 10).Value = 'A'
$ws.Cells.Item(98, 11).Value = 2.05
$ws.Cells.Item(98, 12).Value = 3.2
$ws.Cells.Item(98, 13).Value = 3.5
$ws.Cells.Item(98, 14).Value = 1.75
$ws.Cells.Item(98, 16).Value = 4.2
$ws.Cells.Item(98, 17).Value = -0.5
$ws.Cells.Item(98, 18).Value = 1.8
$ws.Cells.Item(98, 19).Value = 2
$ws.Cells.Item(98, 20).Value = 2.25
$ws.Cells.Item(98, 21).Value = 1.8
$ws.Cells.Item(98, 22).Value = 2
$ws.Cells.Item(98, 23).Value = -1
$ws.Cells.Item(98, 25).Value = 3.2
$ws.Cells.Item(98, 26).Value = -1
$ws.Cells.Item(98, 27).Value = 1
$ws.Cells.Item(98, 28).Value = 0.8
# Row 100
$ws.Cells.Item(100, 2).Value = 6809782
$ws.Cells.Item(100, 6).Value = 'ZulteWaregem'
$ws.Cells.Item(100, 7).Value = 'Anderlecht II'
$ws.Cells.Item(100, 8).Value = 2
$ws.Cells.Item(100, 9).Value = 5
$ws.Cells.Item(100, 11).Value = 1.444
$ws.Cells.Item(100, 12).Value = 4.75
$ws.Cells.Item(100, 13).Value = 5.5
$ws.Cells.Item(100, 14).Value = 1.444
$ws.Cells.Item(100, 15).Value = 4.75
$ws.Cells.Item(100, 16).Value = 5.5
$ws.Cells.Item(100, 17).Value = -1.25
$ws.Cells.Item(100, 18).Value = 1.95
$ws.Cells.Item(100, 19).Value = 1.85
$ws.Cells.Item(100, 20).Value = 3
$ws.Cells.Item(100, 25).Value = 4.5
$ws.Cells.Item(100, 27).Value = 0.8500000000000001
$ws.Cells.Item(100, 28).Value = 0.8
$ws.Cells.Item(100, 29).Value = -1
# Row 101
$ws.Cells.Item(101, 2).Value = 6809785
$ws.Cells.Item(101, 6).Value = 'Francs Borains'
$ws.Cells.Item(101, 7).Value = 'Deinze'
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 11).Value = 4.2
$ws.Cells.Item(101, 12).Value = 3.6
$ws.Cells.Item(101, 13).Value = 1.727
$ws.Cells.Item(101, 14).Value = 4
$ws.Cells.Item(101, 15).Value = 3.5
$ws.Cells.Item(101, 16).Value = 1.8
$ws.Cells.Item(101, 17).Value = 0.5
$ws.Cells.Item(101, 18).Value = 2
$ws.Cells.Item(101, 19).Value = 1.8
$ws.Cells.Item(101, 20).Value = 2.5
$ws.Cells.Item(101, 25).Value = 0.8
$ws.Cells.Item(101, 27).Value = 0.8
$ws.Cells.Item(101, 28).Value = -1
$ws.Cells.Item(101, 29).Value = 1
# Row 108
$ws.Cells.Item(108, 2).Value = 6809791
$ws.Cells.Item(108, 6).Value = 'Deinze'
$ws.Cells.Item(108, 7).Value = 'Club Brugge II'
$ws.Cells.Item(108, 8).Value = 1
$ws.Cells.Item(108, 9).Value = 3
$ws.Cells.Item(108, 11).Value = 1.55
$ws.Cells.Item(108, 12).Value = 4.333
$ws.Cells.Item(108, 13).Value = 4.75
$ws.Cells.Item(108, 14).Value = 1.4
$ws.Cells.Item(108, 15).Value = 4.75
$ws.Cells.Item(108, 16).Value = 6
$ws.Cells.Item(108, 17).Value = -1.25
$ws.Cells.Item(108, 18).Value = 1.875
$ws.Cells.Item(108, 19).Value = 1.975
$ws.Cells.Item(108, 20).Value = 3.25
$ws.Cells.Item(108, 25).Value = 5
$ws.Cells.Item(108, 27).Value = 0.9750000000000001
$ws.Cells.Item(108, 28).Value = 0.9750000000000001
$ws.Cells.Item(108, 29).Value = -1
# Row 109
$ws.Cells.Item(109, 2).Value = 6809788
$ws.Cells.Item(109, 6).Value = 'Anderlecht II'
$ws.Cells.Item(109, 7).Value = 'WaaslandBeveren'
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 9).Value = 2
$ws.Cells.Item(109, 11).Value = 4
$ws.Cells.Item(109, 12).Value = 3.75
$ws.Cells.Item(109, 13).Value = 1.75
$ws.Cells.Item(109, 14).Value = 3.75
$ws.Cells.Item(109, 15).Value = 3.6
$ws.Cells.Item(109, 16).Value = 1.85
$ws.Cells.Item(109, 17).Value = 0.75
$ws.Cells.Item(109, 18).Value = 1.7
$ws.Cells.Item(109, 19).Value = 2.2
$ws.Cells.Item(109, 20).Value = 3
$ws.Cells.Item(109, 25).Value = 0.8500000000000001
$ws.Cells.Item(109, 27).Value = 1.2
$ws.Cells.Item(109, 28).Value = -1
$ws.Cells.Item(109, 29).Value = 0.875
# Row 178
$ws.Cells.Item(178, 2).Value = 6809846
$ws.Cells.Item(178, 6).Value = 'WaaslandBeveren'
$ws.Cells.Item(178, 7).Value = 'Deinze'
$ws.Cells.Item(178, 8).Value = 2
$ws.Cells.Item(178, 9).Value = 2
$ws.Cells.Item(178, 10).Value = 'D'
$ws.Cells.Item(178, 11).Value = 1.909
$ws.Cells.Item(178, 12).Value = 3.5
$ws.Cells.Item(178, 13).Value = 3.5
$ws.Cells.Item(178, 14).Value = 2.2
$ws.Cells.Item(178, 15).Value = 3.4
$ws.Cells.Item(178, 16).Value = 2.9
$ws.Cells.Item(178, 17).Value = -0.25
$ws.Cells.Item(178, 18).Value = 1.975
$ws.Cells.Item(178, 19).Value = 1.825
$ws.Cells.Item(178, 20).Value = 2.75
$ws.Cells.Item(178, 21).Value = 1.95
$ws.Cells.Item(178, 22).Value = 1.85
$ws.Cells.Item(178, 23).Value = -1
$ws.Cells.Item(178, 24).Value = 2.4
$ws.Cells.Item(178, 26).Value = -0.5
$ws.Cells.Item(178, 27).Value = 0.4125
$ws.Cells.Item(178, 28).Value = 0.95
# Row 179
$ws.Cells.Item(179, 2).Value = 6809842
$ws.Cells.Item(179, 6).Value = 'Patro Eisden Maasmechelen'
$ws.Cells.Item(179, 7).Value = 'Lierse Kempenzonen'
$ws.Cells.Item(179, 8).Value = 3
$ws.Cells.Item(179, 9).Value = 0
$ws.Cells.Item(179, 10).Value = 'H'
$ws.Cells.Item(179, 11).Value = 1.533
$ws.Cells.Item(179, 12).Value = 4.2
$ws.Cells.Item(179, 13).Value = 5
$ws.Cells.Item(179, 14).Value = 1.7
$ws.Cells.Item(179, 15).Value = 3.8
$ws.Cells.Item(179, 16).Value = 4.2
$ws.Cells.Item(179, 17).Value = -0.75
$ws.Cells.Item(179, 18).Value = 1.925
$ws.Cells.Item(179, 19).Value = 1.875
$ws.Cells.Item(179, 20).Value = 2.5
$ws.Cells.Item(179, 21).Value = 1.8
$ws.Cells.Item(179, 22).Value = 2
$ws.Cells.Item(179, 23).Value = 0.7
$ws.Cells.Item(179, 24).Value = -1
$ws.Cells.Item(179, 26).Value = 0.925
$ws.Cells.Item(179, 27).Value = -1
$ws.Cells.Item(179, 28).Value = 0.8
# Row 210
$ws.Cells.Item(210, 2).Value = 6809869
$ws.Cells.Item(210, 6).Value = 'KFCO Beerschot Wilrijk'
$ws.Cells.Item(210, 7).Value = 'Genk II'
$ws.Cells.Item(210, 8).Value = 1
$ws.Cells.Item(210, 9).Value = 0
$ws.Cells.Item(210, 10).Value = 'H'
$ws.Cells.Item(210, 11).Value = 1.5
$ws.Cells.Item(210, 13).Value = 6
$ws.Cells.Item(210, 14).Value = 1.333
$ws.Cells.Item(210, 15).Value = 5.25
$ws.Cells.Item(210, 16).Value = 8.5
$ws.Cells.Item(210, 17).Value = -1.5
$ws.Cells.Item(210, 18).Value = 1.85
$ws.Cells.Item(210, 19).Value = 1.95
$ws.Cells.Item(210, 20).Value = 3.5
$ws.Cells.Item(210, 21).Value = 1.975
$ws.Cells.Item(210, 22).Value = 1.825
$ws.Cells.Item(210, 23).Value = 0.333
$ws.Cells.Item(210, 25).Value = -1
$ws.Cells.Item(210, 26).Value = -1
$ws.Cells.Item(210, 27).Value = 0.95
$ws.Cells.Item(210, 29).Value = 0.825
# Row 211
$ws.Cells.Item(211, 2).Value = 6809867
$ws.Cells.Item(211, 6).Value = 'Club Brugge II'
$ws.Cells.Item(211, 7).Value = 'FCV Dender EH'
$ws.Cells.Item(211, 8).Value = 0
$ws.Cells.Item(211, 9).Value = 1
$ws.Cells.Item(211, 10).Value = 'A'
$ws.Cells.Item(211, 11).Value = 6
$ws.Cells.Item(211, 13).Value = 1.5
$ws.Cells.Item(211, 14).Value = 6.5
$ws.Cells.Item(211, 15).Value = 4.5
$ws.Cells.Item(211, 16).Value = 1.5
$ws.Cells.Item(211, 17).Value = 1.25
$ws.Cells.Item(211, 18).Value = 1.8
$ws.Cells.Item(211, 19).Value = 2
$ws.Cells.Item(211, 20).Value = 3
$ws.Cells.Item(211, 21).Value = 1.95
$ws.Cells.Item(211, 22).Value = 1.85
$ws.Cells.Item(211, 23).Value = -1
$ws.Cells.Item(211, 25).Value = 0.5
$ws.Cells.Item(211, 26).Value = 0.4
$ws.Cells.Item(211, 27).Value = -0.5
$ws.Cells.Item(211, 29).Value = 0.8500000000000001
# Row 218
$ws.Cells.Item(218, 2).Value = 6809872
$ws.Cells.Item(218, 6).Value = 'ZulteWaregem'
$ws.Cells.Item(218, 7).Value = 'Patro Eisden Maasmechelen'
$ws.Cells.Item(218, 11).Value = 2
$ws.Cells.Item(218, 12).Value = 3.5
$ws.Cells.Item(218, 13).Value = 3.6
$ws.Cells.Item(218, 14).Value = 2.15
$ws.Cells.Item(218, 15).Value = 3.4
$ws.Cells.Item(218, 16).Value = 3.3
$ws.Cells.Item(218, 17).Value = -0.25
$ws.Cells.Item(218, 18).Value = 1.875
$ws.Cells.Item(218, 19).Value = 1.975
$ws.Cells.Item(218, 20).Value = 2.75
$ws.Cells.Item(218, 21).Value = 1.925
$ws.Cells.Item(218, 22).Value = 1.925
# Row 219
$ws.Cells.Item(219, 2).Value = 6809875
$ws.Cells.Item(219, 6).Value = 'WaaslandBeveren'
$ws.Cells.Item(219, 7).Value = 'Anderlecht II'
$ws.Cells.Item(219, 11).Value = 1.55
$ws.Cells.Item(219, 12).Value = 4
$ws.Cells.Item(219, 13).Value = 5.75
$ws.Cells.Item(219, 14).Value = 1.5
$ws.Cells.Item(219, 15).Value = 4.2
$ws.Cells.Item(219, 16).Value = 6
$ws.Cells.Item(219, 17).Value = -1
$ws.Cells.Item(219, 18).Value = 1.825
$ws.Cells.Item(219, 19).Value = 2.025
$ws.Cells.Item(219, 20).Value = 3
$ws.Cells.Item(219, 21).Value = 1.875
$ws.Cells.Item(219, 22).Value = 1.975
# Row 220
$ws.Cells.Item(220, 21).Value = 1.9
$ws.Cells.Item(220, 22).Value = 1.95
# Row 221
$ws.Cells.Item(221, 14).Value = 1.615
$ws.Cells.Item(221, 15).Value = 4
$ws.Cells.Item(221, 16).Value = 5
# Row 223
$ws.Cells.Item(223, 14).Value = 1.6
$ws.Cells.Item(223, 15).Value = 4.333
$ws.Cells.Item(223, 16).Value = 4.75
$ws.Cells.Item(223, 17).Value = -1
$ws.Cells.Item(223, 21).Value = 1.925
$ws.Cells.Item(223, 22).Value = 1.925
